$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-28 Friday" "2025-11-29 Saturday"

Replace-Text "804÷8=100, 4" "340÷5=68, 0"
Replace-Text "388÷7=55, 3" "376÷5=75, 1"
Replace-Text "196÷8=24, 4" "768÷8=96, 0"
Replace-Text "908÷6=151, 2" "142÷3=47, 1"
Replace-Text "162÷8=20, 2" "787÷7=112, 3"

Replace-Text "854÷3=284, 2" "255÷3=85, 0"
Replace-Text "277÷4=69, 1" "109÷2=54, 1"
Replace-Text "387÷3=129, 0" "475÷6=79, 1"
Replace-Text "878÷7=125, 3" "812÷2=406, 0"
Replace-Text "444÷3=148, 0" "586÷5=117, 1"

Replace-Text "503÷6=83, 5" "439÷5=87, 4"
Replace-Text "127÷3=42, 1" "272÷4=68, 0"
Replace-Text "923÷8=115, 3" "172÷8=21, 4"
Replace-Text "120÷8=15, 0" "484÷7=69, 1"
Replace-Text "842÷6=140, 2" "348÷6=58, 0"

Replace-Text "767÷2=383, 1" "845÷5=169, 0"
Replace-Text "710÷3=236, 2" "873÷3=291, 0"
Replace-Text "539÷9=59, 8" "933÷8=116, 5"
Replace-Text "477÷6=79, 3" "471÷3=157, 0"
Replace-Text "430÷3=143, 1" "763÷3=254, 1"

Replace-Text "839÷4=209, 3" "322÷8=40, 2"
Replace-Text "892÷8=111, 4" "742÷9=82, 4"
Replace-Text "343÷9=38, 1" "157÷6=26, 1"
Replace-Text "398÷7=56, 6" "151÷8=18, 7"
Replace-Text "519÷7=74, 1" "927÷9=103, 0"

$d.Save()
